$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 916 (new weekly price block), shifting existing
# data down by 4 rows (979 -> 983).
$ws.Rows("916:919").Insert()

# New date for this block: 2022-06-02 (Excel serial 44714)
$newDate = Get-Date -Year 2022 -Month 6 -Day 2 -Hour 0 -Minute 0 -Second 0

$qualities = @("Especial", "Primera", "Segunda", "Tercera")
$quantities = @(10, 12, 14, 16)
$prices = @(1875, 1562, 1339, 1172)

for ($i = 0; $i -lt 4; $i++) {
    $r = 916 + $i
    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = $newDate
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108005
    $ws.Cells.Item($r, 10).Value = "Piña"
    $ws.Cells.Item($r, 11).Value = "Caramelo"
    $ws.Cells.Item($r, 12).Value = $qualities[$i]
    $ws.Cells.Item($r, 13).Value = 216
    $ws.Cells.Item($r, 14).Value = 18500
    $ws.Cells.Item($r, 15).Value = 19000
    $ws.Cells.Item($r, 16).Value = 18750
    $ws.Cells.Item($r, 17).Value = "$/caja " + $quantities[$i] + " unidades"
    $ws.Cells.Item($r, 18).Value = "Ecuador"
    $ws.Cells.Item($r, 19).Value = $prices[$i]
    $ws.Cells.Item($r, 20).Value = $quantities[$i]
}

Write-Output "done"
